$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.143.89'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = '1.839.80'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.35'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4635'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3856'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07838'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9609'
$ws.Range('E10').Value = '  -1.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.97'
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '1.842.21'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.674'
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.859'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06904'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.44'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009961'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.69'
$ws.Range('E19').Value = '  -1.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = '28.163.57'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.294'
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.00'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('D25').Value = '2.047.20'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.38'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.14'
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.715'
$ws.Range('E28').Value = '  -5.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.967'
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.69'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09269'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9273'
$ws.Range('E32').Value = '  -4.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.271'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.322'
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.335'
$ws.Range('E35').Value = '  -3.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05795'
$ws.Range('E36').Value = '  -4.64%  '
$ws.Range('E37').Value = '  -4.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.141'
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.749'
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5582'
$ws.Range('E40').Value = '  -1.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.873'
$ws.Range('E41').Value = '  -2.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1759'
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07185'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.59'
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5262'
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.136'
$ws.Range('E46').Value = '  -9.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.122'
$ws.Range('E47').Value = '  -12.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.829'
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '113.67'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.325'
$ws.Range('E51').Value = '  -0.61%  '
